# Sync automatico del tracker (cada 3h): agrega los partidos nuevos
# con sus pronosticos al final de la hoja "Sheet1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ EventId = "14506190"; Fecha = "2025-08-24"; JugadorA = "Luciano Darderi"; JugadorB = "Rinky Hijikata"; Pronostico = "Gana Luciano Darderi"; Cuota = 1.57 },
    @{ EventId = "14506199"; Fecha = "2025-08-24"; JugadorA = "Adrian Mannarino"; JugadorB = "Tallon Griekspoor"; Pronostico = "Gana Tallon Griekspoor"; Cuota = 1.67 },
    @{ EventId = "14506197"; Fecha = "2025-08-24"; JugadorA = "Camilo Ugo Carabelli"; JugadorB = "Tomas Martin Etcheverry"; Pronostico = "Gana Camilo Ugo Carabelli"; Cuota = 3.4 },
    @{ EventId = "14506200"; Fecha = "2025-08-24"; JugadorA = "Jordan Thompson"; JugadorB = "Corentin Moutet"; Pronostico = "Gana Jordan Thompson"; Cuota = 3.4 },
    @{ EventId = "14510039"; Fecha = "2025-08-24"; JugadorA = "Nuria Parrizas Diaz"; JugadorB = "Polina Kudermetova"; Pronostico = "Gana Polina Kudermetova"; Cuota = 1.8 },
    @{ EventId = "14519798"; Fecha = "2025-08-24"; JugadorA = "Oksana Selekhmeteva"; JugadorB = "Marketa Vondrousova"; Pronostico = "Gana Oksana Selekhmeteva"; Cuota = 4 },
    @{ EventId = "14506269"; Fecha = "2025-08-24"; JugadorA = "Jil Teichmann"; JugadorB = "Catherine McNally"; Pronostico = "Gana Jil Teichmann"; Cuota = 3.75 },
    @{ EventId = "14519810"; Fecha = "2025-08-24"; JugadorA = "Rebecca Marino"; JugadorB = "Leylah Fernandez"; Pronostico = "Gana Rebecca Marino"; Cuota = 4.33 },
    @{ EventId = "14519809"; Fecha = "2025-08-24"; JugadorA = "Janice Tjen"; JugadorB = "Veronika Kudermetova"; Pronostico = "Gana Janice Tjen"; Cuota = 2.63 },
    @{ EventId = "14506281"; Fecha = "2025-08-24"; JugadorA = "Lulu Sun"; JugadorB = "Camila Osorio"; Pronostico = "Gana Camila Osorio"; Cuota = 2.1 },
    @{ EventId = "14519799"; Fecha = "2025-08-24"; JugadorA = "Victoria Azarenka"; JugadorB = "Hina Inoue"; Pronostico = "Gana Hina Inoue"; Cuota = 6.5 },
    @{ EventId = "14506275"; Fecha = "2025-08-24"; JugadorA = "Yuliia Starodubtseva"; JugadorB = "Anna Blinkova"; Pronostico = "Gana Yuliia Starodubtseva"; Cuota = 2 }
)

$startRow = 397
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # event_id y fecha se fuerzan a texto (prefijo ') para que no
    # Excel los reinterprete como numero / fecha-serial.
    $ws.Cells.Item($r, 1).Value = "'" + $row.EventId
    $ws.Cells.Item($r, 2).Value = "'" + $row.Fecha
    $ws.Cells.Item($r, 3).Value = $row.JugadorA
    $ws.Cells.Item($r, 4).Value = $row.JugadorB
    $ws.Cells.Item($r, 5).Value = $row.Pronostico
    $ws.Cells.Item($r, 6).Value = $row.Cuota

    # resultado (G) y profit (H) quedan vacios: el partido todavia
    # no se disputo, se completan en una sync posterior.
}
